# Insert a new daily price-log entry for "Piña" at row 266 of the
# "Vega Modelo de Temuco" sheet, pushing the existing rows 266-295 down
# to 267-296 (dimension grows from A1:T295 to A1:T296).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 266..295 down by one to make room for the new record.
$ws.Rows.Item(266).Insert()

# Populate the newly inserted row 266 with the new record's data.
$ws.Range("A266").Value = 10
$ws.Range("B266").Value = "Vega Modelo de Temuco"
$ws.Range("C266").Value = "La Araucanía"
$ws.Range("D266").Value = 44491
$ws.Range("E266").Value = 9
$ws.Range("F266").Value = "Fruta"
$ws.Range("G266").Value = 100108
$ws.Range("H266").Value = "Tropicales y subtropicales"
$ws.Range("I266").Value = 100108005
$ws.Range("J266").Value = "Piña"
$ws.Range("K266").Value = "Caramelo"
$ws.Range("L266").Value = "Primera"
$ws.Range("M266").Value = 200
$ws.Range("N266").Value = 20000
$ws.Range("O266").Value = 20000
$ws.Range("P266").Value = 20000
$ws.Range("Q266").Value = "$/caja 12 unidades"
$ws.Range("R266").Value = "Ecuador"
$ws.Range("S266").Value = 1667
$ws.Range("T266").Value = 12
